# ---- Build the new "2022-Q1" fund-holding sheet (copy of "2021-Q4" layout) ----
$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$srcSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item(6)
$newSheet.Name = "2022-Q1"

# Source sheet has 17 fund rows; target only needs 16, drop the extra trailing row
$newSheet.Rows.Item(18).Delete()

# Fill in fund data rows 2-17
# -- row 2 --
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "007178"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "浙商港股通中华交易服务预期高股息指数增强A"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "7.93"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "90.20"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "8.49"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.6733"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 1

# -- row 3 --
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "007216"
$newSheet.Range("B3").Style = "Normal"
$newSheet.Range("C3").Value = "浙商港股通中华交易服务预期高股息指数增强C"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "4.60"
$newSheet.Range("D3").Style = "Normal"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "90.20"
$newSheet.Range("E3").Style = "Normal"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "8.49"
$newSheet.Range("F3").Style = "Normal"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.3905"
$newSheet.Range("G3").Style = "Normal"
$newSheet.Range("H3").Value = 1

# -- row 4 --
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").NumberFormat = "@"
$newSheet.Range("B4").Value = "007291"
$newSheet.Range("B4").Style = "Normal"
$newSheet.Range("C4").Value = "汇丰晋信港股通双核策略混合"
$newSheet.Range("D4").NumberFormat = "@"
$newSheet.Range("D4").Value = "7.80"
$newSheet.Range("D4").Style = "Normal"
$newSheet.Range("E4").NumberFormat = "@"
$newSheet.Range("E4").Value = "92.66"
$newSheet.Range("E4").Style = "Normal"
$newSheet.Range("F4").NumberFormat = "@"
$newSheet.Range("F4").Value = "4.87"
$newSheet.Range("F4").Style = "Normal"
$newSheet.Range("G4").NumberFormat = "@"
$newSheet.Range("G4").Value = "0.3799"
$newSheet.Range("G4").Style = "Normal"
$newSheet.Range("H4").Value = 6

# -- row 5 --
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").NumberFormat = "@"
$newSheet.Range("B5").Value = "002332"
$newSheet.Range("B5").Style = "Normal"
$newSheet.Range("C5").Value = "汇丰晋信沪港深股票A"
$newSheet.Range("D5").NumberFormat = "@"
$newSheet.Range("D5").Value = "7.90"
$newSheet.Range("D5").Style = "Normal"
$newSheet.Range("E5").NumberFormat = "@"
$newSheet.Range("E5").Value = "92.60"
$newSheet.Range("E5").Style = "Normal"
$newSheet.Range("F5").NumberFormat = "@"
$newSheet.Range("F5").Value = "4.67"
$newSheet.Range("F5").Style = "Normal"
$newSheet.Range("G5").NumberFormat = "@"
$newSheet.Range("G5").Value = "0.3689"
$newSheet.Range("G5").Style = "Normal"
$newSheet.Range("H5").Value = 8

# -- row 6 --
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").NumberFormat = "@"
$newSheet.Range("B6").Value = "202801"
$newSheet.Range("B6").Style = "Normal"
$newSheet.Range("C6").Value = "南方全球精选配置(QDII-FOF)"
$newSheet.Range("D6").NumberFormat = "@"
$newSheet.Range("D6").Value = "18.00"
$newSheet.Range("D6").Style = "Normal"
$newSheet.Range("E6").NumberFormat = "@"
$newSheet.Range("E6").Value = "28.82"
$newSheet.Range("E6").Style = "Normal"
$newSheet.Range("F6").NumberFormat = "@"
$newSheet.Range("F6").Value = "1.04"
$newSheet.Range("F6").Style = "Normal"
$newSheet.Range("G6").NumberFormat = "@"
$newSheet.Range("G6").Value = "0.1872"
$newSheet.Range("G6").Style = "Normal"
$newSheet.Range("H6").Value = 8

# -- row 7 --
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").NumberFormat = "@"
$newSheet.Range("B7").Value = "013009"
$newSheet.Range("B7").Style = "Normal"
$newSheet.Range("C7").Value = "万家港股通精选混合A"
$newSheet.Range("D7").NumberFormat = "@"
$newSheet.Range("D7").Value = "2.97"
$newSheet.Range("D7").Style = "Normal"
$newSheet.Range("E7").NumberFormat = "@"
$newSheet.Range("E7").Value = "81.62"
$newSheet.Range("E7").Style = "Normal"
$newSheet.Range("F7").NumberFormat = "@"
$newSheet.Range("F7").Value = "6.01"
$newSheet.Range("F7").Style = "Normal"
$newSheet.Range("G7").NumberFormat = "@"
$newSheet.Range("G7").Value = "0.1785"
$newSheet.Range("G7").Style = "Normal"
$newSheet.Range("H7").Value = 3

# -- row 8 --
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").NumberFormat = "@"
$newSheet.Range("B8").Value = "001518"
$newSheet.Range("B8").Style = "Normal"
$newSheet.Range("C8").Value = "万家瑞兴灵活配置混合"
$newSheet.Range("D8").NumberFormat = "@"
$newSheet.Range("D8").Value = "3.16"
$newSheet.Range("D8").Style = "Normal"
$newSheet.Range("E8").NumberFormat = "@"
$newSheet.Range("E8").Value = "81.47"
$newSheet.Range("E8").Style = "Normal"
$newSheet.Range("F8").NumberFormat = "@"
$newSheet.Range("F8").Value = "4.47"
$newSheet.Range("F8").Style = "Normal"
$newSheet.Range("G8").NumberFormat = "@"
$newSheet.Range("G8").Value = "0.1413"
$newSheet.Range("G8").Style = "Normal"
$newSheet.Range("H8").Value = 3

# -- row 9 --
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").NumberFormat = "@"
$newSheet.Range("B9").Value = "011355"
$newSheet.Range("B9").Style = "Normal"
$newSheet.Range("C9").Value = "华泰柏瑞港股通时代机遇混合型证券投资基金A"
$newSheet.Range("D9").NumberFormat = "@"
$newSheet.Range("D9").Value = "1.13"
$newSheet.Range("D9").Style = "Normal"
$newSheet.Range("E9").NumberFormat = "@"
$newSheet.Range("E9").Value = "90.93"
$newSheet.Range("E9").Style = "Normal"
$newSheet.Range("F9").NumberFormat = "@"
$newSheet.Range("F9").Value = "8.50"
$newSheet.Range("F9").Style = "Normal"
$newSheet.Range("G9").NumberFormat = "@"
$newSheet.Range("G9").Value = "0.0960"
$newSheet.Range("G9").Style = "Normal"
$newSheet.Range("H9").Value = 1

# -- row 10 --
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").NumberFormat = "@"
$newSheet.Range("B10").Value = "002333"
$newSheet.Range("B10").Style = "Normal"
$newSheet.Range("C10").Value = "汇丰晋信沪港深股票C"
$newSheet.Range("D10").NumberFormat = "@"
$newSheet.Range("D10").Value = "1.23"
$newSheet.Range("D10").Style = "Normal"
$newSheet.Range("E10").NumberFormat = "@"
$newSheet.Range("E10").Value = "92.60"
$newSheet.Range("E10").Style = "Normal"
$newSheet.Range("F10").NumberFormat = "@"
$newSheet.Range("F10").Value = "4.67"
$newSheet.Range("F10").Style = "Normal"
$newSheet.Range("G10").NumberFormat = "@"
$newSheet.Range("G10").Value = "0.0574"
$newSheet.Range("G10").Style = "Normal"
$newSheet.Range("H10").Value = 8

# -- row 11 --
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").NumberFormat = "@"
$newSheet.Range("B11").Value = "013010"
$newSheet.Range("B11").Style = "Normal"
$newSheet.Range("C11").Value = "万家港股通精选混合C"
$newSheet.Range("D11").NumberFormat = "@"
$newSheet.Range("D11").Value = "0.86"
$newSheet.Range("D11").Style = "Normal"
$newSheet.Range("E11").NumberFormat = "@"
$newSheet.Range("E11").Value = "81.62"
$newSheet.Range("E11").Style = "Normal"
$newSheet.Range("F11").NumberFormat = "@"
$newSheet.Range("F11").Value = "6.01"
$newSheet.Range("F11").Style = "Normal"
$newSheet.Range("G11").NumberFormat = "@"
$newSheet.Range("G11").Value = "0.0517"
$newSheet.Range("G11").Style = "Normal"
$newSheet.Range("H11").Value = 3

# -- row 12 --
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").NumberFormat = "@"
$newSheet.Range("B12").Value = "003413"
$newSheet.Range("B12").Style = "Normal"
$newSheet.Range("C12").Value = "华泰柏瑞新经济沪港深灵活配置混合"
$newSheet.Range("D12").NumberFormat = "@"
$newSheet.Range("D12").Value = "0.54"
$newSheet.Range("D12").Style = "Normal"
$newSheet.Range("E12").NumberFormat = "@"
$newSheet.Range("E12").Value = "92.57"
$newSheet.Range("E12").Style = "Normal"
$newSheet.Range("F12").NumberFormat = "@"
$newSheet.Range("F12").Value = "8.84"
$newSheet.Range("F12").Style = "Normal"
$newSheet.Range("G12").NumberFormat = "@"
$newSheet.Range("G12").Value = "0.0477"
$newSheet.Range("G12").Style = "Normal"
$newSheet.Range("H12").Value = 2

# -- row 13 --
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").NumberFormat = "@"
$newSheet.Range("B13").Value = "519139"
$newSheet.Range("B13").Style = "Normal"
$newSheet.Range("C13").Value = "海富通沪港深灵活配置混合"
$newSheet.Range("D13").NumberFormat = "@"
$newSheet.Range("D13").Value = "1.32"
$newSheet.Range("D13").Style = "Normal"
$newSheet.Range("E13").NumberFormat = "@"
$newSheet.Range("E13").Value = "94.37"
$newSheet.Range("E13").Style = "Normal"
$newSheet.Range("F13").NumberFormat = "@"
$newSheet.Range("F13").Value = "3.36"
$newSheet.Range("F13").Style = "Normal"
$newSheet.Range("G13").NumberFormat = "@"
$newSheet.Range("G13").Value = "0.0444"
$newSheet.Range("G13").Style = "Normal"
$newSheet.Range("H13").Value = 8

# -- row 14 --
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").NumberFormat = "@"
$newSheet.Range("B14").Value = "007674"
$newSheet.Range("B14").Style = "Normal"
$newSheet.Range("C14").Value = "工银瑞信产业升级股票A"
$newSheet.Range("D14").NumberFormat = "@"
$newSheet.Range("D14").Value = "0.69"
$newSheet.Range("D14").Style = "Normal"
$newSheet.Range("E14").NumberFormat = "@"
$newSheet.Range("E14").Value = "94.04"
$newSheet.Range("E14").Style = "Normal"
$newSheet.Range("F14").NumberFormat = "@"
$newSheet.Range("F14").Value = "5.43"
$newSheet.Range("F14").Style = "Normal"
$newSheet.Range("G14").NumberFormat = "@"
$newSheet.Range("G14").Value = "0.0375"
$newSheet.Range("G14").Style = "Normal"
$newSheet.Range("H14").Value = 8

# -- row 15 --
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").NumberFormat = "@"
$newSheet.Range("B15").Value = "011356"
$newSheet.Range("B15").Style = "Normal"
$newSheet.Range("C15").Value = "华泰柏瑞港股通时代机遇混合型证券投资基金C"
$newSheet.Range("D15").NumberFormat = "@"
$newSheet.Range("D15").Value = "0.40"
$newSheet.Range("D15").Style = "Normal"
$newSheet.Range("E15").NumberFormat = "@"
$newSheet.Range("E15").Value = "90.93"
$newSheet.Range("E15").Style = "Normal"
$newSheet.Range("F15").NumberFormat = "@"
$newSheet.Range("F15").Value = "8.50"
$newSheet.Range("F15").Style = "Normal"
$newSheet.Range("G15").NumberFormat = "@"
$newSheet.Range("G15").Value = "0.0340"
$newSheet.Range("G15").Style = "Normal"
$newSheet.Range("H15").Value = 1

# -- row 16 --
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").NumberFormat = "@"
$newSheet.Range("B16").Value = "007675"
$newSheet.Range("B16").Style = "Normal"
$newSheet.Range("C16").Value = "工银瑞信产业升级股票C"
$newSheet.Range("D16").NumberFormat = "@"
$newSheet.Range("D16").Value = "0.33"
$newSheet.Range("D16").Style = "Normal"
$newSheet.Range("E16").NumberFormat = "@"
$newSheet.Range("E16").Value = "94.04"
$newSheet.Range("E16").Style = "Normal"
$newSheet.Range("F16").NumberFormat = "@"
$newSheet.Range("F16").Value = "5.43"
$newSheet.Range("F16").Style = "Normal"
$newSheet.Range("G16").NumberFormat = "@"
$newSheet.Range("G16").Value = "0.0179"
$newSheet.Range("G16").Style = "Normal"
$newSheet.Range("H16").Value = 8

# -- row 17 --
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").NumberFormat = "@"
$newSheet.Range("B17").Value = "519602"
$newSheet.Range("B17").Style = "Normal"
$newSheet.Range("C17").Value = "海富通大中华精选混合QDII"
$newSheet.Range("D17").NumberFormat = "@"
$newSheet.Range("D17").Value = "0.11"
$newSheet.Range("D17").Style = "Normal"
$newSheet.Range("E17").NumberFormat = "@"
$newSheet.Range("E17").Value = "89.68"
$newSheet.Range("E17").Style = "Normal"
$newSheet.Range("F17").NumberFormat = "@"
$newSheet.Range("F17").Value = "5.72"
$newSheet.Range("F17").Style = "Normal"
$newSheet.Range("G17").NumberFormat = "@"
$newSheet.Range("G17").Value = "0.0063"
$newSheet.Range("G17").Style = "Normal"
$newSheet.Range("H17").Value = 2

# ---- Update the "总计" (totals) sheet: insert new 2022-Q1 summary row at top ----
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()
$ws.Range("A3:D3").Copy()
$ws.Range("A2:D2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = 16
$ws.Range("D2").Value = 2.71
for ($r = 3; $r -le 7; $r++) {
    $ws.Range("A$r").Value = $r - 2
}
